$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.906.23"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.552.64"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.80"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.485"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.69"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.773.64"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "1.552.66"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "26.905.12"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.47"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.86"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "1.423.86"
$ws.Range("E34").Value = "  +2.94%  "
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.959"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.987"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.63"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.74"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "1.688.33"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.19"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +2.76%  "
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  +6.48%  "
$ws.Range("E51").Value = "  +1.19%  "
